$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 251 - this shifts existing rows 251-265 down to 252-266
$ws.Rows.Item(251).Insert()

# Populate the new row 251 with the new data record
$ws.Cells.Item(251, 1).Value2 = 1
$ws.Cells.Item(251, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(251, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(251, 4).Value2 = 44769
$ws.Cells.Item(251, 5).Value2 = 15
$ws.Cells.Item(251, 6).Value2 = "Fruta"
$ws.Cells.Item(251, 7).Value2 = 100102
$ws.Cells.Item(251, 8).Value2 = "Cítricos"
$ws.Cells.Item(251, 9).Value2 = 100102003
$ws.Cells.Item(251, 10).Value2 = "Limón"
$ws.Cells.Item(251, 11).Value2 = "Sin especificar"
$ws.Cells.Item(251, 12).Value2 = "2a amarillo"
$ws.Cells.Item(251, 13).Value2 = 300
$ws.Cells.Item(251, 14).Value2 = 10000
$ws.Cells.Item(251, 15).Value2 = 11000
$ws.Cells.Item(251, 16).Value2 = 10500
$ws.Cells.Item(251, 17).Value2 = '$/caja 20 kilos'
$ws.Cells.Item(251, 18).Value2 = "Región de Coquimbo"
$ws.Cells.Item(251, 19).Value2 = 525
$ws.Cells.Item(251, 20).Value2 = 20

# Ensure the date-formatted style (same as the other rows in column D) is applied
$ws.Cells.Item(251, 4).NumberFormat = $ws.Cells.Item(252, 4).NumberFormat
